# This workbook's weekly price rows (2-32) got re-shuffled: each row's
# observation data (Fecha, Calidad, Volumen, Precio minimo/maximo/promedio,
# Unidad de comercializacion, Origen, Precio $/Kg, Kg/unidad) is reassigned
# to a different row, while the descriptive columns (A,B,C,E,F,G,H,I,J,K)
# stay constant/unchanged for every row.
#
# Mapping below: destination row -> source row (i.e. "after" row N gets the
# values that "before" row M had).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{
    2  = 21
    3  = 12
    4  = 5
    5  = 30
    6  = 2
    7  = 6
    8  = 22
    9  = 23
    10 = 16
    11 = 20
    12 = 8
    13 = 17
    14 = 18
    15 = 29
    16 = 19
    17 = 26
    18 = 4
    19 = 3
    20 = 24
    21 = 25
    22 = 11
    23 = 7
    24 = 15
    25 = 9
    26 = 10
    27 = 28
    28 = 27
    29 = 32
    30 = 13
    31 = 14
    32 = 31
}

$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot every affected column's value for every row BEFORE any writes,
# since several rows trade values with each other (the map is a permutation,
# not a simple one-directional copy).
$snapshot = @{}
for ($r = 2; $r -le 32; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        if ($c -eq "D") {
            $rowVals[$c] = $ws.Range("$c$r").Value2
        } else {
            $rowVals[$c] = $ws.Range("$c$r").Value()
        }
    }
    $snapshot[$r] = $rowVals
}

# Now write back according to the row map.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
